# Connecting to NedBank Account v4 - update seed Category.xlsx data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "ParentCategoryId" placeholder values in column G for rows 2-10
$ws.Range("G2:G10").ClearContents()

# Populate ParentCategoryId values into column E for rows 11-17 (CategoryTypeId column)
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("E16").Value = 2
$ws.Range("E17").Value = 2

# Clear the old "ParentCategoryId" placeholder values in column G for rows 18-21
$ws.Range("G18:G21").ClearContents()

# Update the active selection/view
$ws.Range("N12").Select() | Out-Null
